# conversie.xlsx - "meer conversie info"
#
# Summary of the change:
#  - A new shared string is introduced: "EindArtikel fallback (wanneer niet
#    met een cijfer begint)".
#  - A new data row is inserted at row 12 (EindArtikel fallback row), pushing
#    the previous rows 12-15 (VerpakkingsArtikel / GrondstofArtikel /
#    PhantomArtikel / ReceptuurArtikel) down to rows 13-16, and the footnote
#    row (B17) down to row 18.
#  - Throughout the data rows (2-16) the "Vergelijk_Begin" (D) and
#    "Vergelijk_ArtikelType" (E) columns are swapped: D now holds the P-code
#    text, E now holds the numeric 69.
#  - Columns D and E are widened to fit their (now textual) header/content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: this sheet's D:E (and most) columns carry a "Text" (@) number
# format, so a plain `.Value = <number>` assignment is stored as text
# (mirrors real Excel: typing a number into a Text-formatted cell keeps it
# textual). To store a genuine numeric value we briefly switch the cell to
# "General", assign, then restore its original number format.
# ---------------------------------------------------------------------
function Set-NumericValue {
    param($range, $value)
    $fmt = $range.NumberFormat
    $range.NumberFormat = "General"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# ---------------------------------------------------------------------
# 1) Insert the new row 12 (shifts old rows 12-17 down by one row).
# ---------------------------------------------------------------------
$ws.Rows.Item(12).Insert()

# ---------------------------------------------------------------------
# 2) Rewrite rows 2-11: swap the D (Vergelijk_Begin / P-code) and
#    E (Vergelijk_ArtikelType / 69) columns.
# ---------------------------------------------------------------------
$pcodes = @("P600","P610","P620","P630","P640","P650","P660","P670","P680","P690")
for ($k = 0; $k -lt 10; $k++) {
    $r = 2 + $k
    $ws.Range("D$r").Value = $pcodes[$k]
    Set-NumericValue $ws.Range("E$r") 69
}

# ---------------------------------------------------------------------
# 3) New row 12: EindArtikel fallback entry.
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "EindArtikel"
$ws.Range("C12").Value = "EindArtikel fallback (wanneer niet met een cijfer begint)"
$ws.Range("D12").Value = "P100"
Set-NumericValue $ws.Range("E12") 69
Set-NumericValue $ws.Range("F12") 6690
Set-NumericValue $ws.Range("G12") 86690
Set-NumericValue $ws.Range("H12") 30669
Set-NumericValue $ws.Range("I12") 76690
Set-NumericValue $ws.Range("J12") 81069

# ---------------------------------------------------------------------
# 4) Rows 13-16 (previously 12-15): fix up the D/E swap, keep existing
#    A/C/F/G/H/I content (Insert already carried it down), and correct the
#    J column to be consistent with row 13's numeric value.
# ---------------------------------------------------------------------
$rowsInfo = @(
    @{ Row = 13; D = "P100"; J = 81020 },
    @{ Row = 14; D = "P100"; J = 81020 },
    @{ Row = 15; D = "P100"; J = $null },
    @{ Row = 16; D = "P100"; J = $null }
)

foreach ($info in $rowsInfo) {
    $r = $info.Row
    $ws.Range("D$r").Value = $info.D
    Set-NumericValue $ws.Range("E$r") 69
    if ($info.J -ne $null) {
        Set-NumericValue $ws.Range("J$r") $info.J
    }
}

# ---------------------------------------------------------------------
# 5) Column widths for D and E (best-fit for the new textual content).
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 27.25
$ws.Columns.Item(5).ColumnWidth = 28.1

# ---------------------------------------------------------------------
# 6) Selection moves to A12 (matches the authored sheetView selection).
# ---------------------------------------------------------------------
$ws.Range("A12").Select()
